# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update India row (row 6) ---
$ws.Range("B6").Value = 2590501
$ws.Range("C6").Value = 1293
$ws.Range("D6").Value = 1862665
$ws.Range("E6").Value = 677737
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = 50099

# --- Update Uzbekistan row (row 61) ---
$ws.Range("B61").Value = 34701
$ws.Range("C61").Value = 173
$ws.Range("D61").Value = 30043
$ws.Range("E61").Value = 4432
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 226

# --- Update Tailandia row (row 117) ---
$ws.Range("B117").Value = 3377
$ws.Range("C117").Value = 1
$ws.Range("D117").Value = 3194

# --- Swap Montserrat / Islas Malvinas ---
# Row 213 previously showed "Islas Malvinas" data, now should show "Montserrat" data
# Row 214 previously showed "Montserrat" data, now should show "Islas Malvinas" data
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# --- Update "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 07:10"
